# Auto commit at 2026-01-11  7:59:58.85
# Updates the "Metrics" sheet's source values (B2:B13) and refreshes the
# selections on "Metrics" and "today" sheets. The "today" sheet's
# B11:B22/E11:E22/F11:F22 cells are formulas that reference Metrics!B2:B13
# (directly or via same-row chains), so they recalculate automatically once
# the Metrics values change - no direct writes needed there.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 192398.17
$metrics.Range("B3").Value  = 144377.70000000001
$metrics.Range("B4").Value  = 46837.490000000005
$metrics.Range("B5").Value  = 7807
$metrics.Range("B6").Value  = 5828268.8999999994
$metrics.Range("B7").Value  = 4915095.33
$metrics.Range("B8").Value  = 1710929.3099999998
$metrics.Range("B9").Value  = 228084
$metrics.Range("B10").Value = 34293649.889999993
$metrics.Range("B11").Value = 32190370.490000002
$metrics.Range("B12").Value = 11992651.35
$metrics.Range("B13").Value = 1325714

# Move the stored selection on Metrics to C9 (matches the target sheetView).
$metrics.Range("C9").Select()

# Move the stored selection on "today" (the active sheet) to F7, and make
# sure it stays the active sheet/tab afterwards.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F7").Select()
